$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column widths (A: 50->49, C: 16->45, add G at width 12) ---
# Note: the ColumnWidth COM property is stored internally as
# (ColumnWidth + 5/6) in the raw OOXML column width, so back off by 5/6 to
# land exactly on the desired stored widths.
$ws.Columns.Item(1).ColumnWidth = 49 - 5/6
$ws.Columns.Item(3).ColumnWidth = 45 - 5/6
$ws.Columns.Item(7).ColumnWidth = 12 - 5/6

# --- New column G: header "Company Name" in G1 (same style as the other
#     header cells, so copy the format from F1 first) ---
$ws.Range("F1").Copy()
$ws.Range("G1").PasteSpecial(-4122)
$ws.Range("G1").Value = "Company Name"

# --- A2: replace billing-address block ---
$ws.Range("A2").Value = "default" + [char]0x2019 + " bill company`n21221 address`ndefault! 1225"

# --- B2: invoice # (must stay text "601", not become the number 601) ---
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "601"
$ws.Range("D2").Copy()
$ws.Range("B2").PasteSpecial(-4122)

# --- C2: new shipping-address block (was empty) ---
$ws.Range("C2").Value = "Orange1 ship company`n1221 address`nOrange1 212"

# --- E2: invoice date (must stay text "01/06/2022", not become a date) ---
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "01/06/2022"
$ws.Range("D2").Copy()
$ws.Range("E2").PasteSpecial(-4122)

# --- F2: total amount (must stay text "448.00", not become the number 448) ---
$ws.Range("F2").NumberFormat = "@"
$ws.Range("F2").Value = "448.00"
$ws.Range("D2").Copy()
$ws.Range("F2").PasteSpecial(-4122)

# --- G2: leave blank (matches the blank D2/empty cell in the same row) ---
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = ""
$ws.Range("D2").Copy()
$ws.Range("G2").PasteSpecial(-4122)
